$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Rebuild the title text as three separate runs, matching how PowerPoint
# splits runs when text is typed incrementally (the middle word "ogg"
# gets flagged by the spell checker):
#   "Audio (mp3, wav, " + "ogg" + ")"
$tr.Text = "Audio (mp3, wav, "
$tr.InsertAfter("ogg") | Out-Null
$tr.InsertAfter(")") | Out-Null
